# Apply feedback corrections to the "Lista de Itens de Trabalho" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lista de Itens de Trabalho")

# Update completion percentages per professor's feedback
$ws.Range("D5").Value = 0.8
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1

# Update the active cell/selection to D5
$ws.Activate()
$ws.Range("D5").Select()
